$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend formatting of the last existing row (99) down to the two new rows
# so the new cells pick up the same date / wrap-text / border styling.
$ws.Range("A99:C99").Copy()
$ws.Range("A100:C101").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New activity rows appended to the calendar table.
# Shared-string authoring order matters: "Ver genero..." (row 101) was
# typed before "Graficos de distribucion..." (row 100), so write the
# text for B101 first, then B100.
$ws.Range("B101").Value = "Ver genero de la base de datos de duque"
$ws.Range("B100").Value = "Graficos de distribucion de los datos demograficos, organizar las graficas de nuevo con los grupos de la base de datos de duque, reunion de avances con el profe"

$ws.Range("A100").Value = Get-Date -Year 2022 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("C100").Value = 9

$ws.Range("A101").Value = Get-Date -Year 2022 -Month 8 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("C101").Value = 1

# Row 100 holds a long wrapped activity description, so it needs a taller
# row (matches the auto-fit height Excel computes for that much wrapped text).
$ws.Rows.Item(100).RowHeight = 55.2

# Update the view state to reflect scrolling to / selecting the new rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D100").Select()
